$d = $word.ActiveDocument

$d.Content.Find.Execute("Expression Specs", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Notation Spec", 2)
